# Apply odds/correct-score updates scraped for 2025-04-04 Flashscore fixtures.
# Each assignment below mirrors one changed <c> cell in the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.5
$ws.Range("I3").Value = 3.2
$ws.Range("V3").Value = 10
$ws.Range("W3").Value = 23
$ws.Range("AE3").Value = 15

# Row 6
$ws.Range("G6").Value = 3.4
$ws.Range("I6").Value = 2.5
$ws.Range("K6").Value = 4.75
$ws.Range("N6").Value = 3.5
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 1.78
$ws.Range("Q6").Value = 2.03
$ws.Range("U6").Value = 15
$ws.Range("Z6").Value = 4.75
$ws.Range("AD6").Value = 5

# Row 8
$ws.Range("G8").Value = 2.8

# Row 10
$ws.Range("G10").Value = 2.35
$ws.Range("I10").Value = 2.75
$ws.Range("Y10").Value = 26
$ws.Range("AA10").Value = 7
$ws.Range("AC10").Value = 41
$ws.Range("AD10").Value = 9.5
$ws.Range("AE10").Value = 15
$ws.Range("AF10").Value = 11

# Row 14
$ws.Range("N14").Value = 2.35
$ws.Range("O14").Value = 1.57

# Row 16
$ws.Range("G16").Value = 2.5
$ws.Range("I16").Value = 2.9
$ws.Range("T16").Value = 8.5
$ws.Range("V16").Value = 10
$ws.Range("X16").Value = 21
$ws.Range("AD16").Value = 9.5

# Row 19
$ws.Range("L19").Value = 1.5
$ws.Range("M19").Value = 2.5
$ws.Range("N19").Value = 2.6
$ws.Range("O19").Value = 1.48

# Row 23
$ws.Range("P23").Value = 1.25
$ws.Range("Q23").Value = 3.75
$ws.Range("Y23").Value = 29
$ws.Range("AA23").Value = 9
$ws.Range("AB23").Value = 13
$ws.Range("AE23").Value = 10
$ws.Range("AI23").Value = 19

# Row 24
$ws.Range("L24").Value = 1.25
$ws.Range("M24").Value = 3.75
$ws.Range("O24").Value = 1.95

# Row 26
$ws.Range("H26").Value = 3.1
$ws.Range("I26").Value = 3.5
$ws.Range("J26").Value = 1.1
$ws.Range("K26").Value = 7
$ws.Range("L26").Value = 1.44
$ws.Range("M26").Value = 2.63
$ws.Range("N26").Value = 2.5
$ws.Range("O26").Value = 1.5
$ws.Range("P26").Value = 1.57
$ws.Range("Q26").Value = 2.25
$ws.Range("R26").Value = 2.05
$ws.Range("S26").Value = 1.7
$ws.Range("T26").Value = 6
$ws.Range("U26").Value = 9.5
$ws.Range("V26").Value = 10
$ws.Range("Y26").Value = 41
$ws.Range("Z26").Value = 7
$ws.Range("AB26").Value = 19
$ws.Range("AC26").Value = 67
$ws.Range("AD26").Value = 8
$ws.Range("AF26").Value = 13
$ws.Range("AG26").Value = 41
$ws.Range("AH26").Value = 34
$ws.Range("AJ26").Value = 501

# Row 28
$ws.Range("J28").Value = 1.08
$ws.Range("K28").Value = 8

# Row 30
$ws.Range("G30").Value = 19
$ws.Range("I30").Value = 1.11
$ws.Range("R30").Value = 2.25
$ws.Range("S30").Value = 1.57
$ws.Range("T30").Value = 45
$ws.Range("U30").Value = 175
$ws.Range("Y30").Value = 175
$ws.Range("Z30").Value = 17
$ws.Range("AA30").Value = 13.5
$ws.Range("AB30").Value = 30
$ws.Range("AD30").Value = 7.6
$ws.Range("AI30").Value = 29

# Row 35
$ws.Range("L35").Value = 1.45
$ws.Range("M35").Value = 2.6
$ws.Range("N35").Value = 2.35
$ws.Range("O35").Value = 1.57
$ws.Range("Q35").Value = 2.5

# Row 39
$ws.Range("L39").Value = 1.22
$ws.Range("M39").Value = 4
$ws.Range("N39").Value = 1.8
$ws.Range("O39").Value = 2

# Row 40
$ws.Range("N40").Value = 2.15
$ws.Range("O40").Value = 1.67
$ws.Range("U40").Value = 13
$ws.Range("AE40").Value = 13

# Row 41
$ws.Range("N41").Value = 1.85
$ws.Range("O41").Value = 2
$ws.Range("T41").Value = 7.5

# Row 42
$ws.Range("AF42").Value = 8.5
$ws.Range("AH42").Value = 13

# Row 44
$ws.Range("J44").Value = 1.1
$ws.Range("K44").Value = 7
$ws.Range("N44").Value = 2.5
$ws.Range("O44").Value = 1.5

# Row 47
$ws.Range("T47").Value = 9.5
$ws.Range("AD47").Value = 8
$ws.Range("AE47").Value = 11

# Row 49
$ws.Range("H49").Value = 3.7
$ws.Range("I49").Value = 4.1
$ws.Range("U49").Value = 9.5
$ws.Range("AB49").Value = 13
$ws.Range("AE49").Value = 21
$ws.Range("AH49").Value = 29

# Row 53
$ws.Range("G53").Value = 2.2
$ws.Range("H53").Value = 2.87
$ws.Range("I53").Value = 3.6
$ws.Range("J53").Value = 1.08
$ws.Range("S53").Value = 2.12
$ws.Range("V53").Value = 8.5
$ws.Range("X53").Value = 17.5
$ws.Range("Y53").Value = 25
$ws.Range("AA53").Value = 5.6
$ws.Range("AB53").Value = 11.5
$ws.Range("AD53").Value = 11
$ws.Range("AI53").Value = 32

# Row 54
$ws.Range("P54").Value = 1.25
$ws.Range("Q54").Value = 3.75
$ws.Range("T54").Value = 23
$ws.Range("AF54").Value = 9
$ws.Range("AG54").Value = 10
$ws.Range("AI54").Value = 21
$ws.Range("AJ54").Value = 151

# Row 56
$ws.Range("G56").Value = 2.15
$ws.Range("H56").Value = 3.2
$ws.Range("I56").Value = 3.2
$ws.Range("L56").Value = 1.4
$ws.Range("M56").Value = 2.5
$ws.Range("T56").Value = 6.2
$ws.Range("U56").Value = 9.25
$ws.Range("W56").Value = 19.5
$ws.Range("X56").Value = 20
$ws.Range("AD56").Value = 7.9
$ws.Range("AE56").Value = 15
$ws.Range("AF56").Value = 12
$ws.Range("AG56").Value = 40
$ws.Range("AH56").Value = 32

# Row 57
$ws.Range("G57").Value = 1.65
$ws.Range("H57").Value = 3.4
$ws.Range("I57").Value = 5.2
$ws.Range("L57").Value = 1.35
$ws.Range("M57").Value = 2.72
$ws.Range("N57").Value = 2.02
$ws.Range("O57").Value = 1.62
$ws.Range("P57").Value = 1.4
$ws.Range("Q57").Value = 2.52
$ws.Range("T57").Value = 5.8
$ws.Range("U57").Value = 7
$ws.Range("W57").Value = 12.5
$ws.Range("X57").Value = 14.5
$ws.Range("Z57").Value = 8.25
$ws.Range("AA57").Value = 6.8
$ws.Range("AB57").Value = 18.5
$ws.Range("AC57").Value = 100
$ws.Range("AD57").Value = 11.75
$ws.Range("AE57").Value = 29
$ws.Range("AF57").Value = 17
$ws.Range("AG57").Value = 100
$ws.Range("AH57").Value = 60
